$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.060.65"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.641.39"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.29"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.25"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  +7.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.399"
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.97"
$ws.Range("E13").Value = "  +5.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.113.88"
$ws.Range("E15").Value = "  +18.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.981.05"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.596.10"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.08"
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("E21").Value = "  +6.11%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.20"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.53"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.15"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0945"
$ws.Range("E30").Value = "  +11.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "519.57"
$ws.Range("E31").Value = "  -8.07%  "
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("E34").Value = "  +7.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.30"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.60"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.26"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.21"
$ws.Range("E42").Value = "  +6.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.03"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.09"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  +5.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.96"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +6.68%  "
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.46"
$ws.Range("E51").Value = "  +2.20%  "
